{"js": "// Remove the sentence about apgdiff/JExamXML test discrepancies and join the\n// remaining \"Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\" + \"u Merlin. \" runs (which\n// are split by the \"_GoBack\" bookmark) into one sentence:\n// \"...apgdiff. Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm serveru Merlin. \"\n\nconst body = context.document.body;\n\n// 1) Delete the whole \"Testy pomoc\u00ed apgdiff ... soubor\u016f. \" sentence that sits\n//    between \"apgdiff. \" and \"Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\".\nconst middleSpan =\n  \"Testy pomoc\u00ed apgdiff vykazovaly rozd\u00edly na m\u00edstech, kde m\u011bl b\u00fdt v\u00fdstup \" +\n  \"spr\u00e1vn\u00fd, pravd\u011bpodobn\u00e1 chyba byla v k\u00f3dov\u00e1n\u00ed porovn\u00e1van\u00fdch soubor\u016f. \";\nconst middleResults = body.search(middleSpan, { matchCase: true });\nmiddleResults.load(\"items\");\nawait context.sync();\n\nif (middleResults.items.length > 0) {\n  middleResults.items[0].delete();\n  await context.sync();\n}\n\n// 2) The paragraph still contains a \"_GoBack\" bookmark between \"...server\"\n//    and the leftover \"u Merlin. \" run. Insert \"u Merlin. \" right before the\n//    bookmark so it joins the preceding \"...server\" text into\n//    \"...serveru Merlin. \".\nconst bookmarkRange = body.getBookmarkRange(\"_GoBack\");\nawait context.sync();\n\nbookmarkRange.insertText(\"u Merlin. \", Word.InsertLocation.before);\nawait context.sync();\n\n// 3) Remove the now-duplicated trailing \"u Merlin. \" run that originally\n//    followed the bookmark.\nconst tailResults = body.search(\"u Merlin. \", { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\n\nif (tailResults.items.length > 0) {\n  tailResults.items[tailResults.items.length - 1].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the sentence about apgdiff/JExamXML test discrepancies and join the\n# remaining \"Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\" + \"u Merlin. \" runs (which\n# are split by the \"_GoBack\" bookmark) into one sentence:\n# \"...apgdiff. Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm serveru Merlin. \"\n\n$d = $word.ActiveDocument\n\n# 1) Duplicate the leftover \"u Merlin. \" text to right before the \"_GoBack\"\n#    bookmark so it is pulled into the same run as \"...na \u0161koln\u00edm server\"\n#    (InsertBefore merges into the immediately preceding run when the\n#    formatting matches), giving \"...na \u0161koln\u00edm serveru Merlin. \" ahead of\n#    the bookmark while the original trailing \"u Merlin. \" run (after the\n#    bookmark) is still there for now.\n$bookmark = $d.Bookmarks.Item(\"_GoBack\")\n$bookmarkStart = $bookmark.Start\n$insertionPoint = $d.Range($bookmarkStart, $bookmarkStart)\n$insertionPoint.InsertBefore(\"u Merlin. \")\n\n# 2) Replace the \"Testy pomoc\u00ed apgdiff ... soubor\u016f. \" sentence (sitting\n#    between \"apgdiff. \" and \"Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\") with\n#    nothing, in the same Find/Replace call that also re-types\n#    \"Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\" so the whole prefix collapses\n#    into a single run.\n$oldSpan = \"Testy pomoc\u00ed apgdiff vykazovaly rozd\u00edly na m\u00edstech, kde m\u011bl b\u00fdt \" + `\n  \"v\u00fdstup spr\u00e1vn\u00fd, pravd\u011bpodobn\u00e1 chyba byla v k\u00f3dov\u00e1n\u00ed porovn\u00e1van\u00fdch \" + `\n  \"soubor\u016f. Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\"\n$newSpan = \"Testov\u00e1n\u00ed prob\u011bhlo na \u0161koln\u00edm server\"\n$range = $d.Content\n$found = $range.Find.Execute($oldSpan, $false, $false, $false, $false, $false, $true, 1, $false, $newSpan, 2)\n\n# 3) Delete the now-redundant original \"u Merlin. \" run that still follows\n#    the bookmark (everything from the bookmark's end to the end of the\n#    document/paragraph).\n$bookmark2 = $d.Bookmarks.Item(\"_GoBack\")\n$tailStart = $bookmark2.End\n$tailEnd = $d.Content.End\n$tailRange = $d.Range($tailStart, $tailEnd)\n$tailRange.Delete()\n"}
